$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'57.889.07"
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Formula = "'3.069.36"
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Range("E3").Value = "  +2.01%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Formula = "'516.96"
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "
$ws.Range("D6").Formula = "'142.50"
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").Formula = "'7.26"
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").Formula = "'3.596.30"
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("E13").Value = "  +3.09%  "
$ws.Range("D14").Formula = "'26.23"
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Formula = "'57.903.77"
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").Formula = "'3.065.62"
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Range("E17").Value = "  +2.16%  "
$ws.Range("D18").Formula = "'6.10"
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Range("E18").Value = "  +2.45%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").Formula = "'331.88"
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Formula = "'0.999"
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Formula = "'0.501"
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").Formula = "'65.59"
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("D25").Formula = "'0.169"
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").Formula = "'0.999"
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Formula = "'0.0₃0909"
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Range("E27").Value = "  -3.50%  "
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Formula = "'7.24"
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Range("E29").Value = "  +5.65%  "
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("E31").Value = "  +3.24%  "
$ws.Range("D32").Formula = "'20.73"
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Range("E32").Value = "  +1.80%  "
$ws.Range("D33").Formula = "'154.78"
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +2.18%  "
$ws.Range("D35").Formula = "'6.00"
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("D36").Formula = "'27.12"
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("E37").Value = "  +3.98%  "
$ws.Range("D38").Formula = "'0.0676"
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("D39").Formula = "'3.111.92"
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("D41").Formula = "'36.61"
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").Formula = "'1.00"
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Formula = "'0.658"
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Formula = "'2.268.52"
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("E45").Value = "  +8.31%  "
$ws.Range("D46").Formula = "'20.97"
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Range("E46").Value = "  +7.68%  "
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").Formula = "'0.943"
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("D49").Formula = "'5.92"
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").Formula = "'0.740"
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Range("E50").Value = "  +9.98%  "
$ws.Range("D51").Formula = "'257.59"
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Range("E51").Value = "  +12.42%  "
